$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 211; this pushes the existing rows 211-306
# down to 212-307 (matching the target diff exactly).
$ws.Rows("211").Insert()

# Populate the newly inserted row 211 with its data.
$ws.Range("A211").Value = 7
$ws.Range("B211").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C211").Value = "Ñuble"
$ws.Range("D211").Value = 45016
$ws.Range("E211").Value = 16
$ws.Range("F211").Value = 100112024
$ws.Range("G211").Value = "Choclo"
$ws.Range("H211").Value = "Choclero"
$ws.Range("I211").Value = "Primera"
$ws.Range("J211").Value = 16000
$ws.Range("K211").Value = 400
$ws.Range("L211").Value = 450
$ws.Range("M211").Value = 425
$ws.Range("N211").Value = "$/unidad"
$ws.Range("O211").Value = "Provincia de Diguillín"
$ws.Range("P211").Value = 425
$ws.Range("Q211").Value = 1
$ws.Range("R211").Value = "Hortaliza"
